# chore: update variable name
# Replace variable name in file transitory_input.xlsx as follows:
#     * MLT_UPPER -> MLT_INCREASE
#     * MLT_LOWER -> MLT_DECREASE
# Row 10 (MLT_INCREASE) also gets its Value filled in (1.2)
# Row 11 (MLT_DECREASE) also gets its Value filled in (0.5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the two variables ---
$ws.Range("A10").Value = "MLT_INCREASE"
$ws.Range("A11").Value = "MLT_DECREASE"

# --- Fill in the previously-empty Value column for those two rows ---
$ws.Range("E10").Value = 1.2
$ws.Range("E11").Value = 0.5

# --- Align the Value column (E7:E13) vertically centered, matching the
#     rest of the sheet's "s=10" (center/center) style now that rows 10/11
#     carry values too ---
$ws.Range("E7:E13").VerticalAlignment = -4108

# --- Update the active selection on the sheet (bottom-right pane) ---
$ws.Range("E19").Select()
